$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.306258201599121
$ws.Range("B1").Value = 1.924432754516602
$ws.Range("C1").Value = 5.152933597564697
$ws.Range("D1").Value = 1.970720291137695
$ws.Range("E1").Value = 1.084469199180603
